$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 356.42856
$ws.Range("I33").Value = 282.5
$ws.Range("K33").Value = 282.5
$ws.Range("M33").Value = -53.5
$ws.Range("H112").Value = 1336.1628
$ws.Range("J112").Value = 1348.9286
$ws.Range("L112").Value = 4046.7858
$ws.Range("N112").Value = -6262.7858
$ws.Range("H118").Value = 982.5
$ws.Range("I118").Value = 810
$ws.Range("J118").Value = 1086
$ws.Range("K118").Value = 2430
$ws.Range("L118").Value = 3258
$ws.Range("M118").Value = -773
$ws.Range("N118").Value = -6572
$ws.Range("H129").Value = 1243.5
$ws.Range("J129").Value = 1299.2262
$ws.Range("L129").Value = 3897.6786
$ws.Range("N129").Value = -13897.6786
$ws.Range("H138").Value = 2970.2415
$ws.Range("I138").Value = 1549.7
$ws.Range("J138").Value = 3717.8948
$ws.Range("K138").Value = 4649.1
$ws.Range("L138").Value = 11153.6844
$ws.Range("M138").Value = 490.8999999999996
$ws.Range("N138").Value = -21433.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4085.1304
$ws.Range("I32").Value = 3998.83
$ws.Range("J32").Value = 4371
$ws.Range("K32").Value = 3998.83
$ws.Range("L32").Value = 4371
$ws.Range("M32").Value = -3711.83
$ws.Range("N32").Value = -4945
$ws.Range("H45").Value = 1916.7858
$ws.Range("I45").Value = 2208.25
$ws.Range("J45").Value = 1528.1666
$ws.Range("K45").Value = 2208.25
$ws.Range("L45").Value = 1528.1666
$ws.Range("M45").Value = -1831.25
$ws.Range("N45").Value = -2282.1666
$ws.Range("H109").Value = 30097.24
$ws.Range("J109").Value = 30097.24
$ws.Range("L109").Value = 30097.24
$ws.Range("N109").Value = -32871.24000000001
$ws.Range("H122").Value = 4021.5356
$ws.Range("I122").Value = 3830.2
$ws.Range("J122").Value = 4499.875
$ws.Range("K122").Value = 11490.6
$ws.Range("L122").Value = 13499.625
$ws.Range("M122").Value = -9040.599999999999
$ws.Range("N122").Value = -18399.625
$ws.Range("H132").Value = 4211.4546
$ws.Range("I132").Value = 3239.125
$ws.Range("J132").Value = 6804.3335
$ws.Range("K132").Value = 9717.375
$ws.Range("L132").Value = 20413.0005
$ws.Range("M132").Value = -7187.375
$ws.Range("N132").Value = -25473.0005
$ws.Range("H137").Value = 38984.285
$ws.Range("J137").Value = 38984.285
$ws.Range("L137").Value = 38984.285
$ws.Range("N137").Value = -49184.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 59000
$ws.Range("J59").Value = 59000
$ws.Range("L59").Value = 59000
$ws.Range("N59").Value = -60694
$ws.Range("H137").Value = 54385.832
$ws.Range("J137").Value = 54385.832
$ws.Range("L137").Value = 54385.832
$ws.Range("N137").Value = -64585.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6537743
$ws.Range("I16").Value = 12347066
$ws.Range("J16").Value = 2254.5
$ws.Range("K16").Value = 12347066
$ws.Range("L16").Value = 2254.5
$ws.Range("M16").Value = -12346779
$ws.Range("N16").Value = -2828.5
$ws.Range("H31").Value = 213359.1
$ws.Range("I31").Value = 588339.3
$ws.Range("J31").Value = 3004.3416
$ws.Range("K31").Value = 588339.3
$ws.Range("L31").Value = 3004.3416
$ws.Range("M31").Value = -588044.3
$ws.Range("N31").Value = -3594.3416
$ws.Range("H34").Value = 213359.1
$ws.Range("I34").Value = 588339.3
$ws.Range("J34").Value = 3004.3416
$ws.Range("K34").Value = 588339.3
$ws.Range("L34").Value = 3004.3416
$ws.Range("M34").Value = -588137.3
$ws.Range("N34").Value = -3408.3416
$ws.Range("H51").Value = 27043.637
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 27043.637
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 27043.637
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -28515.637
$ws.Range("H61").Value = 27043.637
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 27043.637
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 27043.637
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -27739.637
$ws.Range("H64").Value = 43333.332
$ws.Range("J64").Value = 43333.332
$ws.Range("L64").Value = 43333.332
$ws.Range("N64").Value = -43829.332
$ws.Range("H67").Value = 43333.332
$ws.Range("J67").Value = 43333.332
$ws.Range("L67").Value = 43333.332
$ws.Range("N67").Value = -45049.332
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H99").Value = 4683.0713
$ws.Range("I99").Value = 2612
$ws.Range("J99").Value = 6754.143
$ws.Range("K99").Value = 2612
$ws.Range("L99").Value = 6754.143
$ws.Range("M99").Value = -1114
$ws.Range("N99").Value = -9750.143
$ws.Range("H113").Value = 6537743
$ws.Range("I113").Value = 12347066
$ws.Range("J113").Value = 2254.5
$ws.Range("K113").Value = 12347066
$ws.Range("L113").Value = 2254.5
$ws.Range("M113").Value = -12344896
$ws.Range("N113").Value = -6594.5
$ws.Range("H126").Value = 4683.0713
$ws.Range("I126").Value = 2612
$ws.Range("J126").Value = 6754.143
$ws.Range("K126").Value = 7836
$ws.Range("L126").Value = 20262.429
$ws.Range("M126").Value = -5366
$ws.Range("N126").Value = -25202.429
$ws.Range("H134").Value = 9176
$ws.Range("I134").Value = 10120.637
$ws.Range("K134").Value = 30361.911
$ws.Range("M134").Value = -27826.911
$ws.Range("H137").Value = 44485.715
$ws.Range("J137").Value = 44485.715
$ws.Range("L137").Value = 44485.715
$ws.Range("N137").Value = -54685.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3394.94
$ws.Range("I68").Value = 1038.3529
$ws.Range("J68").Value = 4608.9395
$ws.Range("K68").Value = 3115.0587
$ws.Range("L68").Value = 13826.8185
$ws.Range("M68").Value = -2304.0587
$ws.Range("N68").Value = -15448.8185
$ws.Range("H71").Value = 3394.94
$ws.Range("I71").Value = 1038.3529
$ws.Range("J71").Value = 4608.9395
$ws.Range("K71").Value = 9345.176100000001
$ws.Range("L71").Value = 41480.4555
$ws.Range("M71").Value = -5289.176100000001
$ws.Range("N71").Value = -49592.4555
$ws.Range("H131").Value = 858.1900000000001
$ws.Range("J131").Value = 875.0105
$ws.Range("L131").Value = 2625.0315
$ws.Range("N131").Value = -12705.0315
$ws.Range("H137").Value = 2355.5293
$ws.Range("J137").Value = 2272.1428
$ws.Range("L137").Value = 6816.428400000001
$ws.Range("N137").Value = -17016.4284
$ws.Range("H140").Value = 3302.1428
$ws.Range("I140").Value = 846
$ws.Range("J140").Value = 4666.6665
$ws.Range("K140").Value = 2538
$ws.Range("L140").Value = 13999.9995
$ws.Range("M140").Value = 2642
$ws.Range("N140").Value = -24359.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24093.715
$ws.Range("J46").Value = 24004.77
$ws.Range("L46").Value = 24004.77
$ws.Range("N46").Value = -24316.77
$ws.Range("H100").Value = 39775
$ws.Range("J100").Value = 39775
$ws.Range("L100").Value = 39775
$ws.Range("N100").Value = -41939
$ws.Range("H122").Value = 7851.6
$ws.Range("I122").Value = 9999
$ws.Range("J122").Value = 7314.75
$ws.Range("K122").Value = 29997
$ws.Range("L122").Value = 21944.25
$ws.Range("M122").Value = -27547
$ws.Range("N122").Value = -26844.25
$ws.Range("H132").Value = 3860.9
$ws.Range("I132").Value = 2746
$ws.Range("J132").Value = 7205.6
$ws.Range("K132").Value = 8238
$ws.Range("L132").Value = 21616.8
$ws.Range("M132").Value = -5708
$ws.Range("N132").Value = -26676.8
$ws.Range("H137").Value = 30312
$ws.Range("J137").Value = 45780
$ws.Range("L137").Value = 45780
$ws.Range("N137").Value = -55980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
$ws.Range("H100").Value = 2976
$ws.Range("J100").Value = 3268
$ws.Range("L100").Value = 3268
$ws.Range("N100").Value = -4350
$ws.Range("H122").Value = 9149.666999999999
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H136").Value = 3367.6316
$ws.Range("I136").Value = 1499.6875
$ws.Range("K136").Value = 4499.0625
$ws.Range("M136").Value = -1949.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 15666.667
$ws.Range("J63").Value = 15666.667
$ws.Range("L63").Value = 15666.667
$ws.Range("N63").Value = -16914.667
$ws.Range("H66").Value = 15666.667
$ws.Range("J66").Value = 15666.667
$ws.Range("L66").Value = 47000.001
$ws.Range("N66").Value = -53240.001
$ws.Range("H93").Value = 29840.5
$ws.Range("J93").Value = 29840.5
$ws.Range("L93").Value = 29840.5
$ws.Range("N93").Value = -34832.5
$ws.Range("H122").Value = 3492.258
$ws.Range("I122").Value = 2018
$ws.Range("J122").Value = 4706.353
$ws.Range("K122").Value = 6054
$ws.Range("L122").Value = 14119.059
$ws.Range("M122").Value = -3604
$ws.Range("N122").Value = -19019.059
$ws.Range("H132").Value = 4199.625
$ws.Range("I132").Value = 1899.75
$ws.Range("K132").Value = 5699.25
$ws.Range("M132").Value = -3169.25
$ws.Range("H136").Value = 2653.8865
$ws.Range("I136").Value = 1108.5927
$ws.Range("J136").Value = 5108.1763
$ws.Range("K136").Value = 3325.7781
$ws.Range("L136").Value = 15324.5289
$ws.Range("M136").Value = -775.7780999999995
$ws.Range("N136").Value = -20424.5289
